# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers, styled like the rest of the header row.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows 2-50 - season record values are the same for every player on the roster.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 96
    $ws.Cells.Item($row, 31).Value = 66
    $ws.Cells.Item($row, 32).Value = 0
}
